# Generate Report for Archive
# - Update localization "Status" cells from "Ready for handoff" to "In Translation"
#   on the Overview sheet (columns E/F) and on each language sheet (column C).
# - The Status column(s) are narrower since "In Translation" is shorter than
#   "Ready for handoff", so shrink the previously auto-fit column widths to match.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Update the status text wherever it appears.
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsDeDe.Range("C2:C3").Value = "In Translation"

# Shrink the now-narrower Status columns to their new auto-fit width.
$wsOverview.Range("E1").EntireColumn.ColumnWidth = 12.5
$wsOverview.Range("F1").EntireColumn.ColumnWidth = 12.5
$wsZhCn.Range("C1").EntireColumn.ColumnWidth = 12.5
$wsDeDe.Range("C1").EntireColumn.ColumnWidth = 12.5
